# Rename the worksheet to reflect the unified "DataNode" concept
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Adjust row heights (header row gets an explicit height, second header row
# height is reduced slightly)
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 27

# Adjust column widths slightly
$ws.Columns.Item(1).ColumnWidth = 20.08
$ws.Columns.Item(8).ColumnWidth = 25.36

# Update the active selection/cell
$ws.Range("D22").Select() | Out-Null
